$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (no explicit formatting) taken from an untouched data cell,
# used to restore default styling after forcing text number-format on numeric-looking values.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '36.423.07'
$ws.Range("E2").Value = '  -3.24%  '

$ws.Range("D3").Value = '1.981.39'
$ws.Range("E3").Value = '  -4.05%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.01'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  -3.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.628'
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = '  -3.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.14'
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = '  -8.17%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.379'
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = '  -0.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.66'
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = '  -3.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0815'
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = '  +5.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.104'
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = '  -1.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.16'
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = '  +12.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.864'
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = '  -5.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.22'
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = '  -4.88%  '

$ws.Range("D16").Value = '2.271.46'
$ws.Range("E16").Value = '  -4.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.44'
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = '  -2.83%  '

$ws.Range("D18").Value = '1.978.33'
$ws.Range("E18").Value = '  -4.23%  '

$ws.Range("D19").Value = '36.352.41'
$ws.Range("E19").Value = '  -3.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.31'
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = '  -3.67%  '

$ws.Range("D21").Value = '0.0₃0868'
$ws.Range("E21").Value = '  -1.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.31'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  -2.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.13'
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = '  -2.69%  '

$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.61'
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = '  -1.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = '  -2.96%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.88'
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = '  +2.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.46'
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = '  -0.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.134'
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = '  +18.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.88'
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = '  -1.04%  '

$ws.Range("E31").Value = '  -1.22%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.17'
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = '  -3.11%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.94'
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = '  -6.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0627'
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = '  +0.61%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.50'
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = '  -5.29%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.29'
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = '  -8.03%  '

$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.13'
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  -1.42%  '

$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = '  +0.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.79'
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = '  -2.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.10'
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = '  +2.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.24'
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0981'
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  -4.51%  '

$ws.Range("E43").Value = '  -3.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0214'
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = '  -2.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.10'
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = '  -4.81%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.40'
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = '  -4.12%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '92.90'
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  -2.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.64'
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = '  -2.57%  '

$ws.Range("D49").Value = '1.367.90'
$ws.Range("E49").Value = '  -3.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.85'
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  -3.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.18'
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = '  -1.31%  '
